$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text, avoiding Excel's automatic
# number/date inference for numeric-looking strings (e.g. "1.001"),
# then restore the cell style so no stray number-format style is left
# behind (matches the source workbook, which carries no explicit style
# on these data cells).
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "29.959.18"
$ws.Range("E2").Value = "  -0.43%  "
Set-TextValue "D3" "1.870.04"
$ws.Range("E3").Value = "  -2.64%  "
Set-TextValue "D4" "1.001"
Set-TextValue "D5" "319.12"
$ws.Range("E5").Value = "  -3.12%  "
Set-TextValue "D6" "1.000"
$ws.Range("E6").Value = "  +0.06%  "
Set-TextValue "D7" "0.5093"
$ws.Range("E7").Value = "  -2.72%  "
$ws.Range("E8").Value = "  -2.74%  "
Set-TextValue "D9" "0.08202"
$ws.Range("E9").Value = "  -3.41%  "
Set-TextValue "D10" "42.15"
$ws.Range("E10").Value = "  -1.75%  "
Set-TextValue "D11" "1.094"
$ws.Range("E11").Value = "  -3.00%  "
Set-TextValue "D12" "22.95"
$ws.Range("E12").Value = "  +3.02%  "
Set-TextValue "D13" "1.864.35"
$ws.Range("E13").Value = "  -2.95%  "
$ws.Range("E14").Value = "  -1.88%  "
Set-TextValue "D15" "7.193"
$ws.Range("E15").Value = "  -2.83%  "
$ws.Range("E16").Value = "  +0.17%  "
Set-TextValue "D17" "91.99"
$ws.Range("E17").Value = "  -4.53%  "
$ws.Range("E18").Value = "  -2.63%  "
Set-TextValue "D19" "0.06379"
$ws.Range("E19").Value = "  -4.92%  "
Set-TextValue "D20" "17.90"
$ws.Range("E21").Value = "  +0.08%  "
Set-TextValue "D22" "29.944.20"
$ws.Range("E22").Value = "  -0.46%  "
Set-TextValue "D23" "5.827"
$ws.Range("E23").Value = "  -4.02%  "
Set-TextValue "D24" "11.12"
$ws.Range("E24").Value = "  -1.22%  "
Set-TextValue "D25" "2.175"
$ws.Range("E25").Value = "  -2.11%  "
Set-TextValue "D26" "2.088.14"
$ws.Range("E26").Value = "  -2.50%  "
Set-TextValue "D27" "161.25"
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("E28").Value = "  -0.95%  "
Set-TextValue "D29" "2.231"
$ws.Range("E29").Value = "  -9.14%  "
Set-TextValue "D30" "127.47"
$ws.Range("E30").Value = "  -1.53%  "
Set-TextValue "D31" "1.065"
$ws.Range("E31").Value = "  -1.77%  "
$ws.Range("E32").Value = "  -2.46%  "
Set-TextValue "D33" "5.940"
$ws.Range("E33").Value = "  -2.80%  "
Set-TextValue "D34" "3.724"
Set-TextValue "D35" "0.02431"
$ws.Range("E35").Value = "  -3.60%  "
Set-TextValue "D36" "5.205"
$ws.Range("E36").Value = "  -0.49%  "
Set-TextValue "D37" "0.06357"
$ws.Range("E37").Value = "  -3.83%  "
Set-TextValue "D38" "0.2140"
$ws.Range("E38").Value = "  -4.01%  "
Set-TextValue "D39" "1.173"
$ws.Range("E39").Value = "  -5.28%  "
Set-TextValue "D40" "8.494"
$ws.Range("E40").Value = "  -6.03%  "
Set-TextValue "D41" "0.6304"
$ws.Range("E41").Value = "  -4.04%  "
Set-TextValue "D42" "11.28"
$ws.Range("E42").Value = "  -3.74%  "
Set-TextValue "D43" "1.203"
$ws.Range("E43").Value = "  -3.25%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D45" "12.98"
$ws.Range("E45").Value = "  -2.50%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D46" "0.5900"
$ws.Range("E46").Value = "  -4.87%  "
Set-TextValue "D47" "3.635"
$ws.Range("E47").Value = "  -3.69%  "
Set-TextValue "D48" "2.004"
$ws.Range("E48").Value = "  -4.54%  "
Set-TextValue "D49" "122.56"
$ws.Range("E49").Value = "  -2.60%  "
$ws.Range("E50").Value = "  -3.41%  "
Set-TextValue "D51" "1.124"
$ws.Range("E51").Value = "  -2.59%  "
